# "Add Login and Tracking order"
#
# Starting point: workbook has a single sheet "Login" with a small
# username/password table (A1:B3).  This script:
#   1. Adds a new "TrackingOrder" worksheet after "Login" and fills it
#      with an OrderId column.
#   2. Extends the "Login" sheet with more sample username/password rows,
#      several of which look like emails and get turned into mailto:
#      hyperlinks (Excel's usual auto-format for the Hyperlink style).
#
# The exact order of writes below matters: it reproduces the shared
# string table order, and the style (cellXfs) allocation order, seen in
# the target workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # "Login" (already exists)

# --- add the second sheet, right after Login ------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "TrackingOrder"

# --- TrackingOrder header ---------------------------------------------------
$ws2.Range("A1").Value = "OrderId"

# --- Login sheet: extra rows ------------------------------------------------

# Row 2 / col B first (drives shared-string + style allocation order)
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "Hatemyself@1001@@"
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 3
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "thuctanphu12@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:thuctanphu12@gmail.com") | Out-Null

$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "123456"

# Row 4
$ws1.Range("A4").NumberFormat = "@"
$ws1.Range("A4").Value = "thuctanphu12@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A4"), "mailto:thuctanphu12@gmail.com") | Out-Null

$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "Hatemyself@1001@@"
$ws1.Hyperlinks.Add($ws1.Range("B4"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 5
$ws1.Range("A5").Value = "aaa@bbb"
$ws1.Hyperlinks.Add($ws1.Range("A5"), "mailto:aaa@bbb") | Out-Null

$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "Hatemyself@1001@@"
$ws1.Hyperlinks.Add($ws1.Range("B5"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 6 (B6 left empty)
$ws1.Range("A6").NumberFormat = "@"
$ws1.Range("A6").Value = "thuctanphu12@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A6"), "mailto:thuctanphu12@gmail.com") | Out-Null

# Row 2 / col A (written later than rows 3-6 in the real edit history)
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "03547305"

# Row 7 (A7 left empty)
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "Hatemyself@1001@@"
$ws1.Hyperlinks.Add($ws1.Range("B7"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 8
$ws1.Range("A8").Value = 354730579

$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "Hatemyself@1001@@"
$ws1.Hyperlinks.Add($ws1.Range("B8"), "mailto:Hatemyself@1001@@") | Out-Null

# Row 9
$ws1.Range("A9").Value = "thuctanphuaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaaa@gmail.com"

$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "Hatemyself@1001@@"
$ws1.Hyperlinks.Add($ws1.Range("B9"), "mailto:Hatemyself@1001@@") | Out-Null

# --- TrackingOrder sheet: remaining rows -----------------------------------
$ws2.Range("A2").Value = 392921444522425
$ws2.Range("A3").Value = 440277016458018
$ws2.Range("A4").Value = 1
$ws2.Range("A5").Value = "123456789a"

# --- cosmetic touches: column width + selections ---------------------------
$ws1.Columns.Item(1).AutoFit() | Out-Null

$ws1.Range("B9").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("G14").Select() | Out-Null
